$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Copy($ws.Range("A9"))
